$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.678.10"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "'2.284.25"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'304.21"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'95.92"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  -2.61%  "
$ws.Range("D10").Value = "'34.81"
$ws.Range("E10").Value = "  -3.81%  "
$ws.Range("D11").Value = "'0.0782"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "'18.40"
$ws.Range("E12").Value = "  +2.84%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.119"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").Value = "'6.83"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "'2.639.91"
$ws.Range("D16").Value = "'2.287.43"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "'0.773"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "'42.552.34"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "'12.94"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").Value = "'0.0₃0894"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").Value = "'5.97"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("D22").Value = "'67.16"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").Value = "'235.83"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("D27").Value = "'24.70"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").Value = "'2.20"
$ws.Range("E28").Value = "  +7.92%  "
$ws.Range("D29").Value = "'166.25"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").Value = "'8.97"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("D31").Value = "'33.03"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "'17.83"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").Value = "'4.96"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").Value = "'4.48"
$ws.Range("E35").Value = "  -5.39%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").Value = "'0.0683"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").Value = "'1.74"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("D41").Value = "'2.67"
$ws.Range("E41").Value = "  -3.52%  "
$ws.Range("D42").Value = "'1.991.82"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").Value = "'0.0277"
$ws.Range("E43").Value = "  -3.59%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'18.42"
$ws.Range("E44").Value = "  +4.87%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'10.23"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").Value = "'2.07"
$ws.Range("E46").Value = "  -4.63%  "
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("D48").Value = "'2.84"
$ws.Range("E48").Value = "  -4.66%  "
$ws.Range("D49").Value = "'53.58"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").Value = "'2.505.44"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("D51").Value = "'1.13"
$ws.Range("E51").Value = "  +0.66%  "
